# Applies:
#  1. Table style swap on the three tables (slides 14, 15, 16) to
#     {1C1E8AF5-D1D1-4C04-B374-A0CBBD991499}
#  2. Theme colour-scheme update (deck's theme -> standard Office colours)

$p = $ppt.ActivePresentation

function ColorFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1. Table style id updates -------------------------------------------
$newStyleId = "{1C1E8AF5-D1D1-4C04-B374-A0CBBD991499}"
$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colour scheme update ----------------------------------------
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $tcs.Item($i).RGB = ColorFromHex $officeColors[$i - 1]
}
